$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-17 Friday" "2025-10-18 Saturday"

Replace-Text "212×8=1696" "478×4=1912"
Replace-Text "424×9=3816" "892×7=6244"
Replace-Text "616×3=1848" "633×8=5064"
Replace-Text "926×3=2778" "791×7=5537"
Replace-Text "207×4=828" "389×8=3112"

Replace-Text "582×9=5238" "619×7=4333"
Replace-Text "935×8=7480" "222×3=666"
Replace-Text "817×5=4085" "992×4=3968"
Replace-Text "341×2=682" "613×4=2452"
Replace-Text "967×4=3868" "739×5=3695"

Replace-Text "153×2=306" "190×2=380"
Replace-Text "318×3=954" "517×6=3102"
Replace-Text "390×4=1560" "752×6=4512"
Replace-Text "418×8=3344" "284×2=568"
Replace-Text "422×8=3376" "488×7=3416"

Replace-Text "429×7=3003" "437×8=3496"
Replace-Text "186×2=372" "788×3=2364"
Replace-Text "786×4=3144" "227×5=1135"
Replace-Text "189×6=1134" "564×7=3948"
Replace-Text "333×6=1998" "376×3=1128"

Replace-Text "726×8=5808" "114×8=912"
Replace-Text "940×5=4700" "255×8=2040"
Replace-Text "871×5=4355" "407×6=2442"
Replace-Text "919×4=3676" "816×6=4896"
Replace-Text "457×2=914" "841×3=2523"
